$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 561, shifting existing rows 561-639 down to 562-640.
$ws.Rows.Item(561).Insert()

$ws.Cells.Item(561, 1).Value = 3
$ws.Cells.Item(561, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(561, 3).Value = "Coquimbo"
$ws.Cells.Item(561, 4).Value = 44984
$ws.Cells.Item(561, 5).Value = 5
$ws.Cells.Item(561, 6).Value = 100112021
$ws.Cells.Item(561, 7).Value = "Ají"
$ws.Cells.Item(561, 8).Value = "Inferno"
$ws.Cells.Item(561, 9).Value = "Primera"
$ws.Cells.Item(561, 10).Value = 73
$ws.Cells.Item(561, 11).Value = 19000
$ws.Cells.Item(561, 12).Value = 19500
$ws.Cells.Item(561, 13).Value = 19260
$ws.Cells.Item(561, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(561, 15).Value = "Limache"
$ws.Cells.Item(561, 16).Value = 1284
$ws.Cells.Item(561, 17).Value = 15
$ws.Cells.Item(561, 18).Value = "Hortaliza"
